# Updated symbol list on Fri Dec 16 10:53:16 UTC 2022 with GitHub Actions
# Refreshes price/volume-rank text cells in the crypto ticker sheet.
# Price-like ("D" column) values are stored as TEXT (not numbers) in the
# source workbook, so we force text via NumberFormat "@" before writing
# and then restore the "Normal" style so no stray number-format / quote-
# prefix styling is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "250.46"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "24.19"
$ws.Range("D3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.042"
$ws.Range("D4").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.423"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.557"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.327"
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.7967"
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1494"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07906"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03347"
$ws.Range("D12").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09285"
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.563"
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001669"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04765"
$ws.Range("D17").Style = "Normal"

$ws.Range("E18").Value = "17OneONEWorstin24h"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006219"
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005693"
$ws.Range("D20").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.672"
$ws.Range("D23").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1255"
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0006489"
$ws.Range("D27").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04435"
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007051"
$ws.Range("D41").Style = "Normal"

$ws.Range("B42").Value = "BKEXToken"

$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1067"
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"

$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003609"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = "42CEJICEJI"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009222"
$ws.Range("D44").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005892"
$ws.Range("D46").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.7016"
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1167"
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = "48BOLOBOLO"

